# Adds four journals to the research notes worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 66 - Conflict Management and Peace Science
$ws.Range("A66").Value = "Conflict Management and Peace Science"
$ws.Range("B66").Value = "<a href='https://journals.sagepub.com/author-instructions/CMP'target='_blank'>Research Note</a>"
$ws.Range("C66").Value = "6k words"
$ws.Range("D66").Value = 28

# Row 67 - Journal of East Asian Studies
$ws.Range("A67").Value = "Journal of East Asian Studies"
$ws.Range("B67").Value = "<a href='https://www.cambridge.org/core/journals/journal-of-east-asian-studies/information/author-instructions/preparing-your-materials'target='_blank'>Research Note</a>"
$ws.Range("C67").Value = "5k words"
$ws.Range("D67").Value = 13

# Row 68 - International Relations of the Asia-Pacific
$ws.Range("A68").Value = "International Relations of the Asia-Pacific "
$ws.Range("C68").Value = "10k words"
$ws.Range("B68").Value = "<a href='https://academic.oup.com/irap/pages/General_Instructions'target='_blank'>Research Note</a>"
$ws.Range("D68").Value = 18

# Row 69 - Journal of Global Security Studies
$ws.Range("A69").Value = "Journal of Global Security Studies "
$ws.Range("B69").Value = "<a href='https://academic.oup.com/jogss/pages/General_Instructions'target='_blank'>Research Note</a>"
$ws.Range("C69").Value = "3k -- 5k words"
$ws.Range("D69").Value = 28

# Restore the cursor/view state left by the author after the edit
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C72").Select() | Out-Null
$excel.ActiveWindow.Zoom = 170
